# Auto SAP download - update Input sheet values; dependent formulas on
# Template_printout recalc automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Input")

$ws.Range("B3").Value = 20578805
$ws.Range("B4").Value = "A01605"
$ws.Range("B5").Value = "APXCAS2134002"
